# Trade #67 closed at 2026-02-17 08:52:06 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1200.27   # Current Capital
$wsSummary.Range("B4").Value = 0.27      # Total P&L $
$wsSummary.Range("B5").Value = 0.08      # Total P&L %
$wsSummary.Range("B6").Value = 67        # Total Trades
$wsSummary.Range("B7").Value = 28        # Winning Trades
$wsSummary.Range("B9").Value = 41.79     # Win Rate %

# --- Strategy Status sheet (row 4 = MarketMaking) ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 100.27     # Capital
$wsStatus.Range("D4").Value = 67         # Trades
$wsStatus.Range("E4").Value = 0.27       # P&L $
$wsStatus.Range("F4").Value = 0.27       # P&L %
$wsStatus.Range("G4").Value = 41.79      # Win Rate %

# --- All Trades sheet (row 68 = trade #67, now closed) ---
$wsAllTrades = $wb.Worksheets.Item("All Trades")
$wsAllTrades.Range("G68").Value = 0.6
$wsAllTrades.Range("H68").Value = "CLOSED"
$wsAllTrades.Range("I68").Value = 5900
$wsAllTrades.Range("J68").Value = 0.59
$wsAllTrades.Range("K68").Value = 100.27
$wsAllTrades.Range("P68").Value = "early_exit"
$wsAllTrades.Range("Q68").Value = 2.9

# --- MarketMaking sheet (row 68 = trade #67, now closed) ---
$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
$wsMarketMaking.Range("G68").Value = 0.6
$wsMarketMaking.Range("H68").Value = "CLOSED"
$wsMarketMaking.Range("I68").Value = 5900
$wsMarketMaking.Range("J68").Value = 0.59
$wsMarketMaking.Range("K68").Value = 100.27
$wsMarketMaking.Range("P68").Value = "early_exit"
$wsMarketMaking.Range("Q68").Value = 2.9
